# Sheet "Rules": cell C10 ("Integer min" for rule R30) changes from 18 to 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 1
